$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D so "Description" (and "Status") shift right,
# making room for the new "Delivered Amount" column.
$ws.Range("D1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "Delivered Amount"

# Row 2 data
$ws.Range("A2").Value = "CS"
$ws.Range("B2").Value = "Folder"
$ws.Range("C2").Value = 100
$ws.Range("D2").Value = 90
$ws.Range("E2").Value = "Long"
$ws.Range("F2").Value = "Pending"

# Row 3 data
$ws.Range("A3").Value = "CS"
$ws.Range("B3").Value = "Folder"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 99
$ws.Range("E3").Value = "Short"
$ws.Range("F3").Value = "Pending"

# Remove the old row 4 entirely
$ws.Rows.Item(4).Delete()
